$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Fecha" (D) and "Volumen" (M) values between row 3 and row 4
$ws.Range("D3").Value = 44993
$ws.Range("M3").Value = 14

$ws.Range("D4").Value = 44253
$ws.Range("M4").Value = 12
